# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy the existing header formatting (bold,
#     centered, bordered) from A1 onto the three new header cells, then
#     set their text. ---
$ws.Range("A1").Copy()
$ws.Range("AD1").PasteSpecial(-4122)
$ws.Range("AE1").PasteSpecial(-4122)
$ws.Range("AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-44): every player on this roster shares the same team
#     season record, so fill the same W/L/T values down each column. ---
$ws.Range("AD2:AD44").Value = 83
$ws.Range("AE2:AE44").Value = 79
$ws.Range("AF2:AF44").Value = 0
